# DOMA-2542 Localization for Excel template (ticket_report_status_executor)
#
# The "second" ticket row in the analytics template re-used the shared
# strings from the "first" ticket row but with the loop index written as
# "i + 1" (with spaces) instead of "i+1". Normalize the spacing so the
# placeholder matches the expression style used everywhere else in the
# template ({d.tickets[i].*}).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A3").Value = "{d.tickets[i+1].address}"
$ws.Range("B3").Value = "{d.tickets[i+1].processing}"
$ws.Range("C3").Value = "{d.tickets[i+1].completed}"
$ws.Range("D3").Value = "{d.tickets[i+1].canceled}"
$ws.Range("E3").Value = "{d.tickets[i+1].deferred}"
$ws.Range("F3").Value = "{d.tickets[i+1].closed}"
$ws.Range("G3").Value = "{d.tickets[i+1].new_or_reopened}"

# Theme cleanup: the workbook theme's font scheme should keep using
# "Helvetica Neue" as the minor (body) font -- this is what the template's
# default shape/text styles are localized against (replacing the old
# hardcoded "Cambria" references and drop-shadow effect with the theme's
# own minor-font / shadow-less look).
$fontScheme = $wb.Theme.ThemeFontScheme
$fontScheme.MinorFont.Latin = "Helvetica Neue"
